# Apply the edits described in the commit "lagt till tankar om bra lämmelhabitat."
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: nudge the workbook window position (xl/workbook.xml workbookView yWindow 1420 -> 1740).
$excel.ActiveWindow.Top = 1740

# Row 16 (C16): expand on why the metric might be unneeded.
$ws.Range("C16").Value2 = "Behöver jag den här? Kanske onödig. Får se vad NDVI säger först."

# Row 17: rework the lemming-probability description (radius -> "riptriangeln") and add a comment.
$ws.Range("A17").Value2 = "Sannolikhet för lämmel inom riptriangeln runt lyan"
$ws.Range("C17").Value2 = "iptrianglarna eftersom rävar jagar närmare lyan när de har valpar (Frafjord 1993) och måste bära tillbaka mat till lyan (Zapata et al. 1998. Dessutom måste jag hålla observationerna oberoende av varandra. Gallant et al (2014) valde max radius på 1,5 km. "

# Row 18: clarify the threshold wording and add a comment about the threshold value.
$ws.Range("A18").Value2 = "Andel bra lämmelhabitat inom området (sätt gränsvärde för lämmelsannolikhet per pixel)"
$ws.Range("C18").Value2 = "Vad är det maximala sannolikhetsvärdet i en pixel för lämmel under uppgångsfas? I så fall kan jag ta hälften av det som gränsvärde för bra lämmelhabitat."

# Update the selected cell to C18, matching the author's final cursor position.
$ws.Range("C18").Select() | Out-Null
